$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (pushes existing rows 54:112 down to 55:113)
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record
$ws.Cells.Item(54, 1).Value = 9
$ws.Cells.Item(54, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value = 44930
$ws.Cells.Item(54, 5).Value = 13
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100101
$ws.Cells.Item(54, 8).Value = "Berries"
$ws.Cells.Item(54, 9).Value = 100101004
$ws.Cells.Item(54, 10).Value = "Frambuesa"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 300
$ws.Cells.Item(54, 14).Value = 8000
$ws.Cells.Item(54, 15).Value = 8000
$ws.Cells.Item(54, 16).Value = 8000
$ws.Cells.Item(54, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(54, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(54, 19).Value = 4000
$ws.Cells.Item(54, 20).Value = 2
